$d = $word.ActiveDocument
$full = $d.Content
$xml = $full.WordOpenXML
Write-Output ("BEFORE LEN:" + $xml.Length)
Write-Output ("contextualSpacing count before: " + ([regex]::Matches($xml, '<w:contextualSpacing[^/]*/>')).Count)
$fixed = $xml -replace '<w:contextualSpacing[^/]*/>', ''
Write-Output ("contextualSpacing count after: " + ([regex]::Matches($fixed, '<w:contextualSpacing[^/]*/>')).Count)
$full.InsertXML($fixed)
Write-Output "done"
